# TD-937 Fix test data file replacing registrations with enrolments
#
# The "Registrations" column header (table column + header cell B1) is
# renamed to "Enrolments". Renaming the header cell's text is enough for
# Excel to keep the backing ListObject (Table1) column name in sync, so
# there's no need to touch the ListColumns collection separately.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Enrolments"

# Leave the cursor where the author left it when they saved the file.
$ws.Range("B1").Select()
